# "New Test Plans/Try TestCycle"
# A new RMA test cycle ("08SI") was generated by the QA automation tool.
# The "RMA Details Maintenance Grid" sheet's working rows (2-4) are
# re-pointed from the previous cycle's generated RMA / shipper-line /
# record-id values to the new cycle's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 -> RMA-08SI-001 / RMA-08SI-1-1 / a7s5f000000xKZZAA2
$ws.Range("E2").Value = "RMA-08SI-001"
$ws.Range("F2").Value = "RMA-08SI-1-1"
$ws.Range("J2").Value = "a7s5f000000xKZZAA2"

# Row 3 -> RMA-08SI-002 / RMA-08SI-1-2 / a7s5f000000xKZaAAM
$ws.Range("E3").Value = "RMA-08SI-002"
$ws.Range("F3").Value = "RMA-08SI-1-2"
$ws.Range("J3").Value = "a7s5f000000xKZaAAM"

# Row 4 -> RMA-08SI-003 / RMA-08SI-1-3 / a7s5f000000xKZbAAM
$ws.Range("E4").Value = "RMA-08SI-003"
$ws.Range("F4").Value = "RMA-08SI-1-3"
$ws.Range("J4").Value = "a7s5f000000xKZbAAM"

# The "Shipper Line" / "Id" columns are best-fit/auto-sized by the
# authoring tool; the new values re-measure to slightly different
# (sub-character-unit) widths. Nudge the two affected columns to the
# closest widths this engine's character-unit ColumnWidth can express.
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666
$ws.Columns.Item(10).ColumnWidth = 20.5
